$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "97.609.22"
$ws.Range("E2").Value = "  +1.58%  "
$ws.Range("D3").Value = "3.748.73"
$ws.Range("E3").Value = "  +1.95%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "2.46"
$ws.Range("E4").Value = "  +30.70%  "
$ws.Range("E5").Value = "  +0.20%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "233.36"
$ws.Range("E6").Value = "  -0.99%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "662.77"
$ws.Range("E7").Value = "  +2.12%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.454"
$ws.Range("E8").Value = "  +7.86%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "1.18"
$ws.Range("E9").Value = "  +12.91%  "
$ws.Range("E10").Value = "  -0.10%  "
$ws.Range("D11").Value = "3.744.74"
$ws.Range("E11").Value = "  +1.93%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "48.82"
$ws.Range("E12").Value = "  +10.47%  "
$ws.Range("E13").Value = "  +8.79%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.212"
$ws.Range("E14").Value = "  +4.14%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.83"
$ws.Range("E15").Value = "  +2.09%  "
$ws.Range("D16").Value = "4.468.04"
$ws.Range("E16").Value = "  +2.39%  "
$ws.Range("D17").Value = "97.292.19"
$ws.Range("E17").Value = "  +1.40%  "
$ws.Range("E18").Value = "  +15.70%  "
$ws.Range("D19").Value = "3.741.45"
$ws.Range("E19").Value = "  +1.28%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "19.61"
$ws.Range("E20").Value = "  +5.32%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.35"
$ws.Range("E21").Value = "  +3.42%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.559"
$ws.Range("E22").Value = "  +12.16%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "551.61"
$ws.Range("E23").Value = "  +7.18%  "
$ws.Range("B25").Value = "Litecoin"
$ws.Range("C25").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "123.45"
$ws.Range("E25").Value = "  +23.22%  "
$ws.Range("B26").Value = "PEPE"
$ws.Range("C26").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0000219"
$ws.Range("E26").Value = "  +8.31%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.240"
$ws.Range("E27").Value = "  +43.05%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.93"
$ws.Range("E28").Value = "  -0.06%  "
$ws.Range("B29").Value = "Aptos"
$ws.Range("C29").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "13.25"
$ws.Range("E29").Value = "  +0.87%  "
$ws.Range("B30").Value = "InternetComputer(DFINITY)"
$ws.Range("C30").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "13.46"
$ws.Range("E30").Value = "  +12.33%  "
$ws.Range("B31").Value = "PancakeSwap"
$ws.Range("C31").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.09"
$ws.Range("E31").Value = "  +3.22%  "
$ws.Range("B32").Value = "Dai"
$ws.Range("C32").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.00"
$ws.Range("E32").Value = "  +0.04%  "
$ws.Range("B33").Value = "Cronos"
$ws.Range("C33").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.190"
$ws.Range("E33").Value = "  +4.51%  "
$ws.Range("B34").Value = "EthereumClassic"
$ws.Range("C34").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "34.10"
$ws.Range("E34").Value = "  +6.59%  "
$ws.Range("B35").Value = "Fetch.AI"
$ws.Range("C35").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.87"
$ws.Range("E35").Value = "  +2.91%  "
$ws.Range("B36").Value = "PolygonEcosystemToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.638"
$ws.Range("E36").Value = "  +9.54%  "
$ws.Range("B37").Value = "Binance-PegBSC-USD"
$ws.Range("C37").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.997"
$ws.Range("E37").Value = "  -0.58%  "
$ws.Range("B38").Value = "Bittensor"
$ws.Range("C38").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "634.03"
$ws.Range("E38").Value = "  -3.03%  "
$ws.Range("B39").Value = "Filecoin"
$ws.Range("C39").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "7.36"
$ws.Range("E39").Value = "  +9.12%  "
$ws.Range("B40").Value = "RenderToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "8.63"
$ws.Range("E40").Value = "  -1.51%  "
$ws.Range("B41").Value = "USDe"
$ws.Range("C41").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.00"
$ws.Range("E41").Value = "  +0.01%  "
$ws.Range("B42").Value = "Kaspa"
$ws.Range("C42").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.168"
$ws.Range("E42").Value = "  +6.26%  "
$ws.Range("B43").Value = "VeChain"
$ws.Range("C43").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0512"
$ws.Range("E43").Value = "  +15.52%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.07"
$ws.Range("E44").Value = "  +1.72%  "
$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "40.91"
$ws.Range("E45").Value = "  -0.46%  "
$ws.Range("B46").Value = "ARBITRUM"
$ws.Range("C46").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.995"
$ws.Range("E46").Value = "  +4.24%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.486"
$ws.Range("E47").Value = "  +12.81%  "
$ws.Range("B48").Value = "Stacks"
$ws.Range("C48").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.40"
$ws.Range("E48").Value = "  +6.71%  "
$ws.Range("B49").Value = "Cosmos"
$ws.Range("C49").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.17"
$ws.Range("E49").Value = "  +9.45%  "
$ws.Range("B50").Value = "WhiteBITCoin"
$ws.Range("C50").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "23.60"
$ws.Range("E50").Value = "  +0.27%  "
$ws.Range("B51").Value = "OKB"
$ws.Range("C51").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "56.66"
$ws.Range("E51").Value = "  +6.02%  "
